$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 27921
$ws.Range("E2").Value = 1485
$ws.Range("F2").Value = 1466
$ws.Range("G2").Value = 1138
$ws.Range("H2").Value = 636
$ws.Range("I2").Value = 645
$ws.Range("J2").Value = -9
$ws.Range("K2").Value = 20571
$ws.Range("L2").Value = 12270
$ws.Range("M2").Value = 8301
$ws.Range("N2").Value = 8264
$ws.Range("O2").Value = 37
$ws.Range("P2").Value = 500
$ws.Range("Q2").Value = 1618
$ws.Range("R2").Value = -2094
$ws.Range("S2").Value = 659
$ws.Range("T2").Value = 2612
$ws.Range("U2").Value = -995
$ws.Range("V2").Value = 7212
$ws.Range("W2").Value = 5.32
$ws.Range("X2").Value = 2.28
$ws.Range("Y2").Value = 8.039999999999999
$ws.Range("Z2").Value = 3.17
$ws.Range("AA2").Value = 147.82
$ws.Range("AB2").Value = 1602.34
$ws.Range("AC2").Value = 6448
$ws.Range("AD2").Value = 25.12
$ws.Range("AE2").Value = 82805
$ws.Range("AF2").Value = 1.96
$ws.Range("AG2").Value = 1800
$ws.Range("AH2").Value = 1.11
$ws.Range("AI2").Value = 27.94
$ws.Range("AJ2").Value = 8967670

# Row 3
$ws.Range("D3").Value = 26870
$ws.Range("E3").Value = 1553
$ws.Range("F3").Value = 1501
$ws.Range("G3").Value = 1162
$ws.Range("H3").Value = 712
$ws.Range("I3").Value = 718
$ws.Range("J3").Value = -6
$ws.Range("K3").Value = 22518
$ws.Range("L3").Value = 13742
$ws.Range("M3").Value = 8776
$ws.Range("N3").Value = 8745
$ws.Range("O3").Value = 31
$ws.Range("P3").Value = 500
$ws.Range("Q3").Value = 926
$ws.Range("R3").Value = -2307
$ws.Range("S3").Value = 1948
$ws.Range("T3").Value = 2485
$ws.Range("U3").Value = -1559
$ws.Range("V3").Value = 9465
$ws.Range("W3").Value = 5.78
$ws.Range("X3").Value = 2.65
$ws.Range("Y3").Value = 8.44
$ws.Range("Z3").Value = 3.3
$ws.Range("AA3").Value = 156.58
$ws.Range("AB3").Value = 1689.53
$ws.Range("AC3").Value = 7178
$ws.Range("AD3").Value = 20.41
$ws.Range("AE3").Value = 87629
$ws.Range("AF3").Value = 1.67
$ws.Range("AG3").Value = 1800
$ws.Range("AH3").Value = 1.23
$ws.Range("AI3").Value = 25.1
$ws.Range("AJ3").Value = 8967670

# Row 4
$ws.Range("D4").Value = 29283
$ws.Range("E4").Value = 1570
$ws.Range("F4").Value = 1570
$ws.Range("G4").Value = 1152
$ws.Range("H4").Value = 747
$ws.Range("I4").Value = 751
$ws.Range("J4").Value = -4
$ws.Range("K4").Value = 22859
$ws.Range("L4").Value = 13485
$ws.Range("M4").Value = 9374
$ws.Range("N4").Value = 9374
$ws.Range("P4").Value = 500
$ws.Range("Q4").Value = 1522
$ws.Range("R4").Value = -891
$ws.Range("S4").Value = -1405
$ws.Range("T4").Value = 1906
$ws.Range("U4").Value = -384
$ws.Range("V4").Value = 8347
$ws.Range("W4").Value = 5.36
$ws.Range("X4").Value = 2.55
$ws.Range("Y4").Value = 8.289999999999999
$ws.Range("Z4").Value = 3.29
$ws.Range("AA4").Value = 143.87
$ws.Range("AB4").Value = 1818
$ws.Range("AC4").Value = 7514
$ws.Range("AD4").Value = 12.74
$ws.Range("AE4").Value = 93927
$ws.Range("AF4").Value = 1.02
$ws.Range("AG4").Value = 1800
$ws.Range("AH4").Value = 1.88
$ws.Range("AI4").Value = 23.98
$ws.Range("AJ4").Value = 8967670

# Row 5
$ws.Range("D5").Value = 32094
$ws.Range("E5").Value = 1454
$ws.Range("F5").Value = 1454
$ws.Range("G5").Value = 1066
$ws.Range("H5").Value = 679
$ws.Range("I5").Value = 679
$ws.Range("K5").Value = 25824
$ws.Range("L5").Value = 15969
$ws.Range("M5").Value = 9856
$ws.Range("N5").Value = 9856
$ws.Range("P5").Value = 500
$ws.Range("Q5").Value = 274
$ws.Range("R5").Value = -2531
$ws.Range("S5").Value = 2519
$ws.Range("T5").Value = 2185
$ws.Range("U5").Value = -1911
$ws.Range("V5").Value = 10775
$ws.Range("W5").Value = 4.53
$ws.Range("X5").Value = 2.11
$ws.Range("Y5").Value = 7.06
$ws.Range("Z5").Value = 2.79
$ws.Range("AA5").Value = 162.02
$ws.Range("AB5").Value = 1932.69
$ws.Range("AC5").Value = 6786
$ws.Range("AD5").Value = 14.31
$ws.Range("AE5").Value = 98758
$ws.Range("AF5").Value = 0.98
$ws.Range("AG5").Value = 1800
$ws.Range("AH5").Value = 1.85
$ws.Range("AI5").Value = 26.55
$ws.Range("AJ5").Value = 8967670

# Row 6
$ws.Range("D6").Value = 32665
$ws.Range("E6").Value = 704
$ws.Range("F6").Value = 704
$ws.Range("G6").Value = -176
$ws.Range("H6").Value = -531
$ws.Range("I6").Value = -531
$ws.Range("K6").Value = 25097
$ws.Range("L6").Value = 15930
$ws.Range("M6").Value = 9167
$ws.Range("N6").Value = 9167
$ws.Range("P6").Value = 500
$ws.Range("Q6").Value = 1407
$ws.Range("R6").Value = -2330
$ws.Range("S6").Value = 1014
$ws.Range("T6").Value = 2555
$ws.Range("U6").Value = -1148
$ws.Range("V6").Value = 12029
$ws.Range("W6").Value = 2.15
$ws.Range("X6").Value = -1.63
$ws.Range("Y6").Value = -5.59
$ws.Range("Z6").Value = -2.09
$ws.Range("AA6").Value = 173.76
$ws.Range("AB6").Value = 1787.72
$ws.Range("AC6").Value = -5313
$ws.Range("AD6").Value = -10.65
$ws.Range("AE6").Value = 91861
$ws.Range("AF6").Value = 0.62
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 0.44
$ws.Range("AI6").Value = -4.79
$ws.Range("AJ6").Value = 8967670

# Row 7
$ws.Range("D7").Value = 31587
$ws.Range("E7").Value = 850
$ws.Range("G7").Value = 460
$ws.Range("H7").Value = 393
$ws.Range("I7").Value = 393
$ws.Range("K7").Value = 25183
$ws.Range("L7").Value = 15600
$ws.Range("M7").Value = 9583
$ws.Range("N7").Value = 9583
$ws.Range("P7").Value = 500
$ws.Range("Q7").Value = 2570
$ws.Range("R7").Value = -1313
$ws.Range("S7").Value = -350
$ws.Range("T7").Value = 1337
$ws.Range("U7").Value = 1045
$ws.Range("W7").Value = 2.69
$ws.Range("X7").Value = 1.24
$ws.Range("Y7").Value = 4.2
$ws.Range("Z7").Value = 1.57
$ws.Range("AA7").Value = 162.78
$ws.Range("AC7").Value = 3933
$ws.Range("AD7").Value = 12.42
$ws.Range("AE7").Value = 96028
$ws.Range("AF7").Value = 0.51
$ws.Range("AG7").Value = 267
$ws.Range("AH7").Value = 0.55
$ws.Range("AI7").Value = 6.08

# Row 8
$ws.Range("D8").Value = 31477
$ws.Range("E8").Value = 860
$ws.Range("G8").Value = 453
$ws.Range("H8").Value = 340
$ws.Range("I8").Value = 340
$ws.Range("K8").Value = 26143
$ws.Range("L8").Value = 16217
$ws.Range("M8").Value = 9923
$ws.Range("N8").Value = 9923
$ws.Range("P8").Value = 500
$ws.Range("Q8").Value = 2023
$ws.Range("R8").Value = -2307
$ws.Range("S8").Value = -285
$ws.Range("T8").Value = 1270
$ws.Range("U8").Value = 830
$ws.Range("W8").Value = 2.73
$ws.Range("X8").Value = 1.08
$ws.Range("Y8").Value = 3.49
$ws.Range("Z8").Value = 1.32
$ws.Range("AA8").Value = 163.42
$ws.Range("AC8").Value = 3400
$ws.Range("AD8").Value = 14.37
$ws.Range("AE8").Value = 99435
$ws.Range("AF8").Value = 0.49
$ws.Range("AG8").Value = 300
$ws.Range("AH8").Value = 0.61
$ws.Range("AI8").Value = 7.91

# Row 9
$ws.Range("D9").Value = 31793
$ws.Range("E9").Value = 910
$ws.Range("G9").Value = 497
$ws.Range("H9").Value = 377
$ws.Range("I9").Value = 377
$ws.Range("K9").Value = 26137
$ws.Range("L9").Value = 15833
$ws.Range("M9").Value = 10303
$ws.Range("N9").Value = 10303
$ws.Range("P9").Value = 500
$ws.Range("Q9").Value = 2303
$ws.Range("R9").Value = -1500
$ws.Range("S9").Value = -285
$ws.Range("T9").Value = 1480
$ws.Range("U9").Value = 1005
$ws.Range("W9").Value = 2.86
$ws.Range("X9").Value = 1.19
$ws.Range("Y9").Value = 3.72
$ws.Range("Z9").Value = 1.44
$ws.Range("AA9").Value = 153.67
$ws.Range("AC9").Value = 3767
$ws.Range("AD9").Value = 12.97
$ws.Range("AE9").Value = 103243
$ws.Range("AF9").Value = 0.47
$ws.Range("AG9").Value = 300
$ws.Range("AH9").Value = 0.61
$ws.Range("AI9").Value = 7.14

# Clear cells that no longer exist in the updated data
$ws.Range("O4").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()